$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 data: Date, Start Time, End Time (formula in D6 recalculates automatically)
$ws.Range("A6").Value = 43052
$ws.Range("B6").Value = 0.92708333333333337
$ws.Range("C6").Value = 0.95833333333333337
$ws.Range("B6:C6").NumberFormat = "h:mm AM/PM"

# Move the active selection to F9 (matches the saved sheetView selection)
$ws.Range("F9").Select()
